$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 260 (weekly update); this shifts the
# previous rows 260-274 down to 261-275, preserving the row that
# previously fell off the end (old row 274) as the new row 275.
$ws.Rows.Item(260).Insert()

# Populate the newly inserted row 260 with this week's record.
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44516
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112024
$ws.Cells.Item(260, 7).Value = "Choclo"
$ws.Cells.Item(260, 8).Value = "Dulce o Americano"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 55
$ws.Cells.Item(260, 11).Value = 32000
$ws.Cells.Item(260, 12).Value = 32000
$ws.Cells.Item(260, 13).Value = 32000
$ws.Cells.Item(260, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(260, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(260, 16).Value = 457
$ws.Cells.Item(260, 17).Value = 70
$ws.Cells.Item(260, 18).Value = "Hortaliza"
